$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a duplicated header layout:
#   Row 1 = generic "ColumnN" placeholders (from the original query refresh)
#   Row 2 = the real lowercase column headers (hotel_name, locality, price, rating, score, review)
#   Rows 3-82 = the 80 data rows
#
# The edit removes the redundant lowercase header row (row 2), shifting the
# 80 data rows up to become rows 2-81, and replaces row 1 with proper
# Title Case column headers that match the table's column names.

# Delete row 2 (the lowercase duplicate header row); rows below shift up.
$ws.Rows("2:2").Delete()

# Re-label the header row (row 1) with the proper Title Case names,
# matching the table column order used for the Booking_com table.
$ws.Range("F1").Value = "No of Reviews"
$ws.Range("G1").Value = "Reviews"
$ws.Range("E1").Value = "Score"
$ws.Range("D1").Value = "Rating"
$ws.Range("C1").Value = "Price"
$ws.Range("B1").Value = "Locality"
$ws.Range("A1").Value = "Hotel_name"

# Fix up the workbook-level defined name that tracked the query table's
# external data range so it reflects the new (smaller) extent.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet2!ExternalData_1") {
        $n.RefersTo = "=Sheet2!`$A`$1:`$G`$81"
    }
}

# Restore the active selection to the cell the author ended up on.
$ws.Range("C9").Select()
